$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(152).Insert()

$ws.Cells.Item(152, 1).Value = 10
$ws.Cells.Item(152, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(152, 3).Value = "La Araucanía"
$ws.Cells.Item(152, 4).Value = 44841
$ws.Cells.Item(152, 5).Value = 9
$ws.Cells.Item(152, 6).Value = 100114007
$ws.Cells.Item(152, 7).Value = "Jengibre"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 50
$ws.Cells.Item(152, 11).Value = 20000
$ws.Cells.Item(152, 12).Value = 20000
$ws.Cells.Item(152, 13).Value = 20000
$ws.Cells.Item(152, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(152, 15).Value = "Perú"
$ws.Cells.Item(152, 16).Value = 1538
$ws.Cells.Item(152, 17).Value = 13
$ws.Cells.Item(152, 18).Value = "Hortaliza"
